$wb = $excel.ActiveWorkbook

# Sheets where column B should simply mirror column A (value + style),
# for every data row below the header row.
$mirrorSheets = @{
    "Bristol"       = 40
    "Leeds"         = 8
    "LondonCentral" = 5
    "LondonEast"    = 5
    "LondonSouth"   = 4
    "Manchester"    = 16
    "MidlandsEast"  = 15
    "MidlandsWest"  = 10
    "Newcastle"     = 10
    "Wales"         = 66
    "Watford"       = 14
}

foreach ($name in $mirrorSheets.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $lastRow = $mirrorSheets[$name]
    $ws.Range("A2:A$lastRow").Copy($ws.Range("B2:B$lastRow"))
}

# Scotland sheet gets real managing-office addresses in column B instead of
# a straight copy of column A. Assigned bottom-to-top so the new shared
# strings land in the same order as the target workbook.
$scotland = $wb.Worksheets.Item("Scotland")
$scotland.Range("B5").Value = "54-56 Melville Street, Edinburgh, EH3 7HF"
$scotland.Range("B4").Value = "Ground Floor, Block C, Caledonian House, Greenmarket, Dundee, DD1 4QB"
$scotland.Range("B3").Value = "Ground Floor, AB1, 48 Huntly Street, Aberdeen, AB10 1SH"
$scotland.Range("B2").Value = "Eagle Building, 215 Bothwell Street, Glasgow, G2 7TS"

# Reset the selection on every sheet to C1 (matches the saved file state).
foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Select() | Out-Null
}

# Make Scotland the active sheet/tab, matching the target workbook state.
$scotland.Activate() | Out-Null
$scotland.Range("C1").Select() | Out-Null
